$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row (row 1): "_old" suffix -> "_FV2410",
#    "_new" suffix -> "_FV2504" (the "diff" header is unchanged).
# ---------------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the used range into an Excel Table ("Table1") so the header row
#    gets the table/autofilter structure + column names.
# ---------------------------------------------------------------------------
$dataRange = $ws.UsedRange
$tbl = $ws.ListObjects.Add(1, $dataRange, 0, 1)                 # xlSrcRange, xlYes
$tbl.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split after row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
